$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 1.045634197174252
$ws.Range("C2").Value = 0.275542105873825
$ws.Range("D2").Value = 0.03046775130748358
$ws.Range("E2").Value = 0.1124537920805
$ws.Range("F2").Value = 0.7163731155385165
$ws.Range("I2").Value = 0.6502193167962744
$ws.Range("L2").Value = 0.2034208097478185
$ws.Range("M2").Value = 0.2252972918433827
$ws.Range("O2").Value = 2.464290182500434
$ws.Range("B3").Value = 0.9366361524037643
$ws.Range("C3").Value = 0.2575815738631491
$ws.Range("D3").Value = 0.02857376532227107
$ws.Range("E3").Value = 0.1136927504555906
$ws.Range("F3").Value = 0.7173328932413128
$ws.Range("I3").Value = 0.6606873729560725
$ws.Range("L3").Value = 0.2008182403035406
$ws.Range("M3").Value = 0.2083546590957894
$ws.Range("O3").Value = 2.482262635119852
$ws.Range("B4").Value = 0.8696382720132192
$ws.Range("C4").Value = 0.2464920295393256
$ws.Range("D4").Value = 0.02740426704868071
$ws.Range("E4").Value = 0.1144972361491424
$ws.Range("F4").Value = 0.7184973022593297
$ws.Range("I4").Value = 0.6675917603935879
$ws.Range("L4").Value = 0.1993230431602839
$ws.Range("M4").Value = 0.1979837443195791
$ws.Range("O4").Value = 2.495242534327303
$ws.Range("B5").Value = 0.8423195618109389
$ws.Range("C5").Value = 0.2419577229974834
$ws.Range("D5").Value = 0.02692606250356988
$ws.Range("E5").Value = 0.1148360884150527
$ws.Range("F5").Value = 0.7191162298728386
$ws.Range("I5").Value = 0.6705251319483203
$ws.Range("L5").Value = 0.1987396368325562
$ws.Range("M5").Value = 0.1937658041536707
$ws.Range("O5").Value = 2.501020252061565
$ws.Range("B6").Value = 0.8377823638057293
$ws.Range("C6").Value = 0.2412038937217744
$ws.Range("D6").Value = 0.02684655980055339
$ws.Range("E6").Value = 0.1148930204348664
$ws.Range("F6").Value = 0.7192277193780114
$ws.Range("I6").Value = 0.6710194457080902
$ws.Range("L6").Value = 0.1986443283784638
$ws.Range("M6").Value = 0.1930659269799335
$ws.Range("O6").Value = 2.502009108945728
$ws.Range("B7").Value = 0.8692699069517857
$ws.Range("C7").Value = 0.2464309395641067
$ws.Range("D7").Value = 0.02739782435102711
$ws.Range("E7").Value = 0.1145017613998934
$ws.Range("F7").Value = 0.7185050648335292
$ws.Range("I7").Value = 0.6676308361035446
$ws.Range("L7").Value = 0.199315070207625
$ws.Range("M7").Value = 0.1979268257804065
$ws.Range("O7").Value = 2.495318478485046
$ws.Range("B8").Value = 1.008067751678652
$ws.Range("C8").Value = 0.2693622887165077
$ws.Range("D8").Value = 0.02981608920237733
$ws.Range("E8").Value = 0.1128719109717631
$ws.Range("F8").Value = 0.7165845669065618
$ws.Range("I8").Value = 0.653729616116415
$ws.Range("L8").Value = 0.2025021387519388
$ws.Range("M8").Value = 0.2194489965804109
$ws.Range("O8").Value = 2.470083019476675
$ws.Range("B9").Value = 1.279610971120462
$ws.Range("C9").Value = 0.3138304408695376
$ws.Range("D9").Value = 0.03450493096153906
$ws.Range("E9").Value = 0.1100223379822185
$ws.Range("F9").Value = 0.7173913841943005
$ws.Range("I9").Value = 0.630260642262364
$ws.Range("L9").Value = 0.2095659013155498
$ws.Range("M9").Value = 0.2618979411096802
$ws.Range("O9").Value = 2.436061116283099
$ws.Range("B10").Value = 1.478660359083108
$ws.Range("C10").Value = 0.3461852634428055
$ws.Range("D10").Value = 0.03791609509145388
$ws.Range("E10").Value = 0.1081390733635226
$ws.Range("F10").Value = 0.7207869119655257
$ws.Range("I10").Value = 0.615338292680061
$ws.Range("L10").Value = 0.2152504419322554
$ws.Range("M10").Value = 0.2932246992110024
$ws.Range("O10").Value = 2.420542698247374
$ws.Range("B11").Value = 1.569102509158256
$ws.Range("C11").Value = 0.3608336005262629
$ws.Range("D11").Value = 0.03946036856835633
$ws.Range("E11").Value = 0.1073277917433326
$ws.Range("F11").Value = 0.7229435143738883
$ws.Range("I11").Value = 0.6090554886343504
$ws.Range("L11").Value = 0.2179436861985806
$ws.Range("M11").Value = 0.3075046140601927
$ws.Range("O11").Value = 2.41555173266633
$ws.Range("B12").Value = 1.603333855341759
$ws.Range("C12").Value = 0.3663702309458756
$ws.Range("D12").Value = 0.04004404255829286
$ws.Range("E12").Value = 0.1070270982521404
$ws.Range("F12").Value = 0.7238483971279948
$ws.Range("I12").Value = 0.6067492080171277
$ws.Range("L12").Value = 0.2189789422516668
$ws.Range("M12").Value = 0.3129160324503175
$ws.Range("O12").Value = 2.413960055005589
$ws.Range("B13").Value = 1.595962313461996
$ws.Range("C13").Value = 0.3651782848843936
$ws.Range("D13").Value = 0.03991838775980483
$ws.Range("E13").Value = 0.1070915681119611
$ws.Range("F13").Value = 0.7236495870778299
$ws.Range("I13").Value = 0.6072426622645466
$ws.Range("L13").Value = 0.2187552979637388
$ws.Range("M13").Value = 0.3117504165152951
$ws.Range("O13").Value = 2.414289571086215
$ws.Range("B14").Value = 1.571919096420288
$ws.Range("C14").Value = 0.3612893121557477
$ws.Range("D14").Value = 0.03950841018987461
$ws.Range("E14").Value = 0.107302922874819
$ws.Range("F14").Value = 0.7230161901957217
$ws.Range("I14").Value = 0.6088642877646713
$ws.Range("L14").Value = 0.2180285493234493
$ws.Range("M14").Value = 0.3079497376315743
$ws.Range("O14").Value = 2.415414801287824
$ws.Range("B15").Value = 1.557189646414656
$ws.Range("C15").Value = 0.3589058436480457
$ws.Range("D15").Value = 0.03925714171094086
$ws.Range("E15").Value = 0.1074332327105663
$ws.Range("F15").Value = 0.722639712360639
$ws.Range("I15").Value = 0.6098670773259265
$ws.Range("L15").Value = 0.2175853966418231
$ws.Range("M15").Value = 0.305622214989512
$ws.Range("O15").Value = 2.416142909029674
$ws.Range("B16").Value = 1.472747453689919
$ws.Range("C16").Value = 0.3452265246717445
$ws.Range("D16").Value = 0.03781501980694912
$ws.Range("E16").Value = 0.1081930056654929
$ws.Range("F16").Value = 0.7206583064177181
$ws.Range("I16").Value = 0.6157590765512531
$ws.Range("L16").Value = 0.2150765873835638
$ws.Range("M16").Value = 0.2922920373933664
$ws.Range("O16").Value = 2.420910564972672
$ws.Range("B17").Value = 1.420916340722727
$ws.Range("C17").Value = 0.3368165494905782
$ws.Range("D17").Value = 0.03692838484893457
$ws.Range("E17").Value = 0.1086707284572508
$ws.Range("F17").Value = 0.7195996765182784
$ws.Range("I17").Value = 0.6195032256755582
$ws.Range("L17").Value = 0.2135649658057019
$ws.Range("M17").Value = 0.2841216905597932
$ws.Range("O17").Value = 2.424365753637574
$ws.Range("B18").Value = 1.39109454764872
$ws.Range("C18").Value = 0.3319727803993828
$ws.Range("D18").Value = 0.03641771356466705
$ws.Range("E18").Value = 0.1089497782171074
$ws.Range("F18").Value = 0.7190483678259056
$ws.Range("I18").Value = 0.6217043245301817
$ws.Range("L18").Value = 0.2127056252984971
$ws.Range("M18").Value = 0.279425089466443
$ws.Range("O18").Value = 2.426547727489066
$ws.Range("B19").Value = 1.380995755283436
$ws.Range("C19").Value = 0.3303316437249748
$ws.Range("D19").Value = 0.03624468942756209
$ws.Range("E19").Value = 0.1090449944510394
$ws.Range("F19").Value = 0.7188715879432621
$ws.Range("I19").Value = 0.6224577417810444
$ws.Range("L19").Value = 0.212416404286401
$ws.Range("M19").Value = 0.2778353846108885
$ws.Range("O19").Value = 2.427319909992775
$ws.Range("B20").Value = 1.426434889908307
$ws.Range("C20").Value = 0.337712488296404
$ws.Range("D20").Value = 0.03702284157967028
$ws.Range("E20").Value = 0.1086194315536604
$ws.Range("F20").Value = 0.7197064077778847
$ws.Range("I20").Value = 0.6190997305470098
$ws.Range("L20").Value = 0.2137248351101988
$ws.Range("M20").Value = 0.2849911536357439
$ws.Range("O20").Value = 2.423977792118563
$ws.Range("B21").Value = 1.578981651589174
$ws.Range("C21").Value = 0.3624318817693393
$ws.Range("D21").Value = 0.03962886089509965
$ws.Range("E21").Value = 0.1072406659582483
$ws.Range("F21").Value = 0.7231998380498652
$ws.Range("I21").Value = 0.6083859975027366
$ws.Range("L21").Value = 0.218241595920702
$ws.Range("M21").Value = 0.3090659848330475
$ws.Range("O21").Value = 2.415076191547712
$ws.Range("B22").Value = 1.678579074148104
$ws.Range("C22").Value = 0.3785268095889762
$ws.Range("D22").Value = 0.04132556677707555
$ws.Range("E22").Value = 0.1063775688458906
$ws.Range("F22").Value = 0.7259973379648983
$ws.Range("I22").Value = 0.6018088509698885
$ws.Range("L22").Value = 0.2212832007076031
$ws.Range("M22").Value = 0.3248230270383274
$ws.Range("O22").Value = 2.410997551451743
$ws.Range("B23").Value = 1.625431870907221
$ws.Range("C23").Value = 0.3699422913196599
$ws.Range("D23").Value = 0.04042060650529322
$ws.Range("E23").Value = 0.1068347460142145
$ws.Range("F23").Value = 0.7244571235114847
$ws.Range("I23").Value = 0.6052802543888305
$ws.Range("L23").Value = 0.2196516531520558
$ws.Range("M23").Value = 0.3164112064830888
$ws.Range("O23").Value = 2.413014992727113
$ws.Range("B24").Value = 1.423940026560786
$ws.Range("C24").Value = 0.3373074616351062
$ws.Range("D24").Value = 0.03698014059573751
$ws.Range("E24").Value = 0.1086426091755779
$ws.Range("F24").Value = 0.719657976069783
$ws.Range("I24").Value = 0.6192819995162608
$ws.Range("L24").Value = 0.2136525279464507
$ws.Range("M24").Value = 0.2845980673648185
$ws.Range("O24").Value = 2.424152580487174
$ws.Range("B25").Value = 1.206226402676577
$ws.Range("C25").Value = 0.3018553224590903
$ws.Range("D25").Value = 0.03324231280945611
$ws.Range("E25").Value = 0.1107562205917422
$ws.Range("F25").Value = 0.7166820004923764
$ws.Range("I25").Value = 0.6362028030439113
$ws.Range("L25").Value = 0.2075679747205257
$ws.Range("M25").Value = 0.2503892155922003
$ws.Range("O25").Value = 2.44360417470611
